$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with new values
$ws.Range("A1").Value = "Linh Đàm"
$ws.Range("A2").Value = "Thanh Trì"

# Add new rows
$ws.Range("A3").Value = "Hoàng Mai"
$ws.Range("A4").Value = "Hà Đông"

# Update font size of A1 to 12 (creates new style)
$ws.Range("A1").Font.Size = 12
$ws.Rows.Item(1).RowHeight = 15.75

# Update selection to A4 (last edited cell)
$ws.Range("A4").Select()
